$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) holds plain text, not numbers (e.g. "23.894.36",
# "1.002", "308.50"). Writing a numeric-looking string straight into a
# General-formatted cell makes Excel parse it as a number, which both
# re-types the cell AND can silently drop a significant trailing zero
# (308.50 -> 308.5). A leading quote forces literal text, just like a
# user typing an apostrophe before the value, so every price below is
# entered that way to keep it verbatim text.

$ws.Range("D2").Value = "'23.894.36"
$ws.Range("E2").Value = "  +0.23%  "
$ws.Range("D3").Value = "'1.648.11"
$ws.Range("E3").Value = "  +2.07%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "'308.50"
$ws.Range("E5").Value = "  +0.47%  "
$ws.Range("D6").Value = "'1.002"
$ws.Range("E6").Value = "  +0.15%  "
$ws.Range("D7").Value = "'0.3882"
$ws.Range("E7").Value = "  -0.48%  "
$ws.Range("D8").Value = "'0.3829"
$ws.Range("E8").Value = "  +0.37%  "
$ws.Range("D9").Value = "'50.99"
$ws.Range("E9").Value = "  +4.13%  "
$ws.Range("E10").Value = "  -0.69%  "
$ws.Range("D11").Value = "'1.002"
$ws.Range("E11").Value = "  +0.17%  "
$ws.Range("D12").Value = "'0.08450"
$ws.Range("E12").Value = "  +0.36%  "
$ws.Range("D13").Value = "'23.86"
$ws.Range("E13").Value = "  -0.06%  "
$ws.Range("D14").Value = "'7.100"
$ws.Range("E14").Value = "  +1.25%  "
$ws.Range("D15").Value = "'7.790"
$ws.Range("E15").Value = "  +3.67%  "
$ws.Range("D16").Value = "'0.00001312"
$ws.Range("E16").Value = "  +2.99%  "
$ws.Range("D17").Value = "'1.647.44"
$ws.Range("E17").Value = "  +2.25%  "
$ws.Range("D18").Value = "'94.38"
$ws.Range("E18").Value = "  +1.46%  "
$ws.Range("D19").Value = "'0.06989"
$ws.Range("E19").Value = "  +1.19%  "
$ws.Range("D20").Value = "'19.70"
$ws.Range("E20").Value = "  -1.58%  "
$ws.Range("D21").Value = "'6.865"
$ws.Range("E21").Value = "  +1.01%  "
$ws.Range("D22").Value = "'1.002"
$ws.Range("E22").Value = "  +0.24%  "
$ws.Range("D23").Value = "'13.57"
$ws.Range("E23").Value = "  +1.28%  "
$ws.Range("D24").Value = "'23.894.97"
$ws.Range("E24").Value = "  +0.19%  "
$ws.Range("D25").Value = "'2.481"
$ws.Range("E25").Value = "  +2.00%  "
$ws.Range("D26").Value = "'3.038"
$ws.Range("E26").Value = "  +7.26%  "
$ws.Range("D27").Value = "'22.11"
$ws.Range("E27").Value = "  +0.03%  "
$ws.Range("E28").Value = "  -2.96%  "
$ws.Range("D29").Value = "'5.432"
$ws.Range("E29").Value = "  +3.50%  "
$ws.Range("D30").Value = "'139.21"
$ws.Range("E30").Value = "  +0.18%  "
$ws.Range("D31").Value = "'7.754"
$ws.Range("E31").Value = "  -0.91%  "
$ws.Range("D32").Value = "'2.495"
$ws.Range("E32").Value = "  +0.29%  "
$ws.Range("D33").Value = "'1.826.63"
$ws.Range("E33").Value = "  +2.00%  "
$ws.Range("D34").Value = "'1.024"
$ws.Range("E34").Value = "  +5.67%  "
$ws.Range("D35").Value = "'0.08024"
$ws.Range("E35").Value = "  -0.50%  "
$ws.Range("D36").Value = "'0.02951"
$ws.Range("E36").Value = "  +2.71%  "
$ws.Range("D37").Value = "'6.680"
$ws.Range("E37").Value = "  +1.52%  "
$ws.Range("D38").Value = "'10.90"
$ws.Range("E38").Value = "  +5.43%  "
$ws.Range("D39").Value = "'0.2672"
$ws.Range("E39").Value = "  +0.68%  "
$ws.Range("D40").Value = "'0.09109"
$ws.Range("E40").Value = "  -0.67%  "
$ws.Range("D43").Value = "'1.419"
$ws.Range("E43").Value = "  -0.46%  "
$ws.Range("D44").Value = "'16.25"
$ws.Range("E44").Value = "  +2.93%  "
$ws.Range("D45").Value = "'0.6912"
$ws.Range("E45").Value = "  +1.27%  "
$ws.Range("D46").Value = "'2.444"
$ws.Range("E46").Value = "  +0.05%  "
$ws.Range("D47").Value = "'4.070"
$ws.Range("E47").Value = "  +0.34%  "
$ws.Range("D48").Value = "'1.001"
$ws.Range("E48").Value = "  +0.10%  "
$ws.Range("D49").Value = "'0.08269"
$ws.Range("E49").Value = "  +0.46%  "
$ws.Range("D50").Value = "'133.97"
$ws.Range("E50").Value = "  +0.82%  "
$ws.Range("D51").Value = "'1.222"
$ws.Range("E51").Value = "  +1.81%  "

# Row 41/42 swap: Aptos <-> TheSandbox, with refreshed price/volume data
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "'0.7531"
$ws.Range("E41").Value = "  +1.28%  "

$ws.Range("B42").Value = "Aptos"
$ws.Range("C42").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D42").Value = "'13.42"
$ws.Range("E42").Value = "  -1.07%  "
